# option-zeh.test-cases.xlsx — add "zeh02" test-case sheet + matching
# rows on the "answers" sheet, mirroring the existing "zeh01" pattern.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Append 13 new rows (15-27) to the "answers" sheet, describing the
#    new "zeh02" test case. Formatting is copied from the analogous
#    "zeh01" rows (2 = first row, 3:13 = body rows, 14 = last row) so the
#    borders/styles match exactly, then values/text are written on top.
# ---------------------------------------------------------------------
$answers = $wb.Worksheets.Item("answers")

$answers.Range("A2:F2").Copy()
$answers.Range("A15:F15").PasteSpecial(-4122)

$answers.Range("A3:F13").Copy()
$answers.Range("A16:F26").PasteSpecial(-4122)

$answers.Range("A14:F14").Copy()
$answers.Range("A27:F27").PasteSpecial(-4122)

$answers.Range("A15").Value = "zeh02"
$answers.Range("B15").Value = "housingAnswer"
$answers.Range("C15").Value = "residentCount"
$answers.Range("D15").Value = 5
$answers.Range("E15").Value = "number"

$answers.Range("A16").Value = "zeh02"
$answers.Range("B16").Value = "housingAnswer"
$answers.Range("C16").Value = "housingSizeKey"
$answers.Range("D16").Value = "4-room"
$answers.Range("E16").Value = "string"

$answers.Range("A17").Value = "zeh02"
$answers.Range("B17").Value = "housingAnswer"
$answers.Range("C17").Value = "electricityIntensityKey"
$answers.Range("D17").Value = "30-renewable"
$answers.Range("E17").Value = "string"

$answers.Range("A18").Value = "zeh02"
$answers.Range("B18").Value = "housingAnswer"
$answers.Range("C18").Value = "electricityMonthlyConsumption"
$answers.Range("D18").Value = 750
$answers.Range("E18").Value = "number"

$answers.Range("A19").Value = "zeh02"
$answers.Range("B19").Value = "housingAnswer"
$answers.Range("C19").Value = "electricitySeasonFactorKey"
$answers.Range("D19").Value = "march"
$answers.Range("E19").Value = "string"

$answers.Range("A20").Value = "zeh02"
$answers.Range("B20").Value = "housingAnswer"
$answers.Range("C20").Value = "useGas"
$answers.Range("D20").Value = $true
$answers.Range("E20").Value = "boolean"

$answers.Range("A21").Value = "zeh02"
$answers.Range("B21").Value = "housingAnswer"
$answers.Range("C21").Value = "energyHeatIntensityKey"
$answers.Range("D21").Value = "lpg"
$answers.Range("E21").Value = "string"

$answers.Range("A22").Value = "zeh02"
$answers.Range("B22").Value = "housingAnswer"
$answers.Range("C22").Value = "gasMonthlyConsumption"
$answers.Range("D22").Value = 15
$answers.Range("E22").Value = "number"

$answers.Range("A23").Value = "zeh02"
$answers.Range("B23").Value = "housingAnswer"
$answers.Range("C23").Value = "gasSeasonFactorKey"
$answers.Range("D23").Value = "december"
$answers.Range("E23").Value = "string"

$answers.Range("A24").Value = "zeh02"
$answers.Range("B24").Value = "housingAnswer"
$answers.Range("C24").Value = "useKerosene"
$answers.Range("D24").Value = $false
$answers.Range("E24").Value = "boolean"

$answers.Range("A25").Value = "zeh02"
$answers.Range("B25").Value = "housingAnswer"
$answers.Range("C25").Value = "keroseneMonthlyConsumption"
$answers.Range("D25").Value = 200
$answers.Range("E25").Value = "number"

$answers.Range("A26").Value = "zeh02"
$answers.Range("B26").Value = "housingAnswer"
$answers.Range("C26").Value = "keroseneMonthCount"
$answers.Range("D26").Value = 2
$answers.Range("E26").Value = "number"

$answers.Range("A27").Value = "zeh02"
$answers.Range("B27").Value = "housingAnswer"
$answers.Range("C27").Value = "housingAmountByRegionFirstKey"
$answers.Range("D27").Value = "northeast"
$answers.Range("E27").Value = "string"

# ---------------------------------------------------------------------
# 2) Duplicate the "zeh01" worksheet to create "zeh02" (placed right
#    after "zeh01", matching the sheet order in the workbook).
# ---------------------------------------------------------------------
$zeh01 = $wb.Worksheets.Item("zeh01")
$zeh01.Copy($null, $zeh01)
$zeh02 = $wb.Worksheets.Item($wb.Worksheets.Count)
$zeh02.Name = "zeh02"

# ---------------------------------------------------------------------
# 3) Restore/update the selections on each sheet. Selecting a range
#    makes that sheet active, so "zeh02" (the sheet that should end up
#    active/selected) is selected last.
# ---------------------------------------------------------------------
$answers.Range("C21").Select()
$zeh01.Range("M17").Select()
$zeh02.Range("G16").Select()
